# Update PSM-DID regression and parallel trends test results (v3.0)
# Re-run after DID variable reconstruction: updates coefficient, std_error,
# t_stat, p_value, ci_lower, ci_upper for each event-time row, and the
# significance marker for the -3 and >=+5 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function SetRow($Row, $Coefficient, $StdError, $TStat, $PValue, $CiLower, $CiUpper) {
    $arr = New-Object 'object[,]' 1,6
    $arr[0,0] = $Coefficient
    $arr[0,1] = $StdError
    $arr[0,2] = $TStat
    $arr[0,3] = $PValue
    $arr[0,4] = $CiLower
    $arr[0,5] = $CiUpper
    $rng = "D" + $Row + ":I" + $Row
    $ws.Range($rng).Value = $arr
}

# relative_year = -6 (<= -5)
SetRow 2 -0.001368844122222731 0.04753860046202624 -0.02879437149850806 0.9770307038803834 -0.09454450102779414 0.09180681278334869

# relative_year = -4
SetRow 3 -0.005108385660693202 0.03332914192935848 -0.1532708424213407 0.8781960710041736 -0.07043350384223582 0.06021673252084941

# relative_year = -3
SetRow 4 -0.05990535617960974 0.02783971199299977 -2.151795111769576 0.03150132862582344 -0.1144711916858893 -0.005339520673330191
$ws.Range("J4").Value = "**"

# relative_year = -2
SetRow 5 -0.01744856034482041 0.01629785451482956 -1.070604743031901 0.2844422491512688 -0.04939235519388635 0.01449523450424553

# relative_year = -1 (baseline period) -- row 6 unchanged

# relative_year = 0
SetRow 7 -0.01465799438367035 0.01191238916555551 -1.23048317008092 0.2186227044640374 -0.03800627714815914 0.008690288380818444

# relative_year = 1
SetRow 8 -0.002259670361465111 0.0180400021528455 -0.1252588742684095 0.9003278727300974 -0.0376180745810423 0.03309873385811207

# relative_year = 2
SetRow 9 0.02723364902570055 0.02282674303673263 1.193058903842583 0.2329504376623222 -0.01750676732629541 0.07197406537769652

# relative_year = 3
SetRow 10 0.02369062793671469 0.0245753719337353 0.963998754549627 0.3351323079475659 -0.0244771010534065 0.07185835692683587

# relative_year = 4
SetRow 11 0.005717451143604203 0.02776009823010392 0.2059593268082899 0.8368381075376876 -0.04869234138739948 0.06012724367460789

# relative_year = 6 (>= +5)
SetRow 12 0.04567381451318967 0.03654167614995177 1.249910221024439 0.2114399542900685 -0.0259478707407158 0.1172954997670951
$ws.Range("J12").ClearContents()
